$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Lot Track" value in B2 to the new "Pro-SYDATA1 (Lot track)" label
$ws.Range("B2").Value = "Pro-SYDATA1 (Lot track)"

# Move the active selection from A3:XFD3 (row select) to just B3
$null = $ws.Range("B3").Select()
